# Update priors and targets with latest data
# Applies the cell-value edits on the "constant" sheet (rows 41 and 43)
# and updates the active selection, matching the target revision.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("constant")
$ws.Activate()

# Row 41 (prev_se_subclin_lowinf_cxr): prior/target values updated
$ws.Range("B41").Value = 0.25
$ws.Range("D41").Value = 0.07
$ws.Range("E41").Value = 0.52

# Row 43 (prev_se_subclin_inf_cxr): prior/target values updated
$ws.Range("B43").Value = 0.85
$ws.Range("D43").Value = 0.72
$ws.Range("E43").Value = 0.93

# Match the saved cursor/selection position from the edit
$ws.Range("B45").Select()
